$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the crypto price / 1h-volume figures (and the RenderToken <->
# WEMIXToken / Aave <-> Maker row swaps) to match the latest scrape.
#
# Every new value is written with a leading apostrophe so the COM layer
# stores it as literal text -- exactly like the original inlineStr cells --
# instead of re-interpreting numeric-looking text (e.g. "245.54",
# "37.041.92") as a number.
$ws.Range("D2").Value = '''37.041.92'
$ws.Range("E2").Value = '''  +1.20%  '
$ws.Range("D3").Value = '''1.981.95'
$ws.Range("E3").Value = '''  +1.12%  '
$ws.Range("E4").Value = '''  +0.05%  '
$ws.Range("D5").Value = '''245.54'
$ws.Range("E5").Value = '''  +0.48%  '
$ws.Range("E6").Value = '''  +2.06%  '
$ws.Range("D7").Value = '''61.35'
$ws.Range("E7").Value = '''  +4.64%  '
$ws.Range("E8").Value = '''  -0.05%  '
$ws.Range("D9").Value = '''0.382'
$ws.Range("E9").Value = '''  +1.05%  '
$ws.Range("D10").Value = '''0.0798'
$ws.Range("E10").Value = '''  -1.23%  '
$ws.Range("E11").Value = '''  -0.10%  '
$ws.Range("D12").Value = '''14.94'
$ws.Range("E12").Value = '''  +9.20%  '
$ws.Range("D13").Value = '''22.37'
$ws.Range("E13").Value = '''  +0.81%  '
$ws.Range("D14").Value = '''0.842'
$ws.Range("E14").Value = '''  +2.01%  '
$ws.Range("D15").Value = '''2.274.48'
$ws.Range("E15").Value = '''  +1.14%  '
$ws.Range("D16").Value = '''5.45'
$ws.Range("E16").Value = '''  +3.52%  '
$ws.Range("D17").Value = '''1.982.90'
$ws.Range("E17").Value = '''  +1.16%  '
$ws.Range("D18").Value = '''36.918.20'
$ws.Range("E18").Value = '''  +1.17%  '
$ws.Range("E19").Value = '''  +0.44%  '
$ws.Range("D20").Value = '''0.0₃0860'
$ws.Range("E20").Value = '''  +0.39%  '
$ws.Range("D21").Value = '''5.16'
$ws.Range("E21").Value = '''  +2.13%  '
$ws.Range("D22").Value = '''230.26'
$ws.Range("E22").Value = '''  +0.83%  '
$ws.Range("E23").Value = '''  +0.02%  '
$ws.Range("D24").Value = '''2.50'
$ws.Range("E24").Value = '''  +1.82%  '
$ws.Range("E25").Value = '''  +0.57%  '
$ws.Range("D26").Value = '''0.152'
$ws.Range("E26").Value = '''  +10.06%  '
$ws.Range("D27").Value = '''9.24'
$ws.Range("E27").Value = '''  +0.03%  '
$ws.Range("D28").Value = '''163.09'
$ws.Range("E28").Value = '''  +1.88%  '
$ws.Range("D29").Value = '''19.61'
$ws.Range("E29").Value = '''  +1.00%  '
$ws.Range("D30").Value = '''1.35'
$ws.Range("E30").Value = '''  +15.64%  '
$ws.Range("E31").Value = '''  +1.80%  '
$ws.Range("D32").Value = '''4.85'
$ws.Range("E32").Value = '''  +3.01%  '
$ws.Range("E33").Value = '''  +0.58%  '
$ws.Range("D34").Value = '''4.51'
$ws.Range("E34").Value = '''  +5.19%  '
$ws.Range("D35").Value = '''2.28'
$ws.Range("E35").Value = '''  +1.78%  '
$ws.Range("E36").Value = '''  +0.00%  '
$ws.Range("B37").Value = '''WEMIXToken'
$ws.Range("C37").Value = '''https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D37").Value = '''1.78'
$ws.Range("E37").Value = '''  +0.60%  '
$ws.Range("B38").Value = '''RenderToken'
$ws.Range("C38").Value = '''https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").Value = '''3.30'
$ws.Range("E38").Value = '''  -1.53%  '
$ws.Range("E39").Value = '''  -6.24%  '
$ws.Range("D40").Value = '''0.0973'
$ws.Range("E40").Value = '''  -0.88%  '
$ws.Range("E41").Value = '''  +1.13%  '
$ws.Range("D42").Value = '''1.17'
$ws.Range("E42").Value = '''  +0.12%  '
$ws.Range("E43").Value = '''  +0.77%  '
$ws.Range("D44").Value = '''16.52'
$ws.Range("E44").Value = '''  +3.17%  '
$ws.Range("B45").Value = '''Maker'
$ws.Range("C45").Value = '''https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value = '''1.372.23'
$ws.Range("E45").Value = '''  +0.80%  '
$ws.Range("B46").Value = '''Aave'
$ws.Range("C46").Value = '''https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = '''90.10'
$ws.Range("E46").Value = '''  +2.73%  '
$ws.Range("D47").Value = '''1.03'
$ws.Range("E47").Value = '''  -0.07%  '
$ws.Range("D48").Value = '''7.23'
$ws.Range("E48").Value = '''  +1.41%  '
$ws.Range("E49").Value = '''  -0.39%  '
$ws.Range("D50").Value = '''46.25'
$ws.Range("E50").Value = '''  +6.12%  '
$ws.Range("D51").Value = '''1.94'
$ws.Range("E51").Value = '''  +8.58%  '

# The apostrophe trick marks the cells with a "quote prefix" text style;
# clear that back off so the cells keep their original (unstyled) look.
$ws.Range("B2:E51").ClearFormats()
